{"js": "// Populate the \"Kit Components\" table with the correct reference/overview\n// values, replacing the placeholder content that had been written into the\n// wrong table, and blanking out the cells that no longer apply.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Locate the \"Kit Components\" table (header row: Description / Quantity /\n// Volume / Storage of opened/reconstituted material) rather than assuming a\n// fixed index.\nconst headerCells = tables.items.map((t) => t.getCell(0, 0));\nheaderCells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nlet kitTable = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  if (headerCells[i].value === \"Description\") {\n    kitTable = tables.items[i];\n    break;\n  }\n}\n\nconst newValues = [\n  [\"Reactive Species\", \"Mouse\", \"\", \"\"],\n  [\"Size\", \"96 wells/kit, with removable strips.\", \"\", \"\"],\n  [\"Sensitivity*\", \"12 pg/ml\", \"\", \"\"],\n  [\n    \"Storage Instructions\",\n    \"Store at 4\u00b0C for 6 months, at -20\u00b0C for 12 months. Avoid multiple freeze-thaw cycles (Ships with gel ice, can store for up to 3 days in room temperature. Freeze upon receiving.)\",\n    \"\",\n    \"\",\n  ],\n  [\"Uniprot ID\", \"P15947\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\"],\n];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = kitTable.getCell(r + 1, c);\n    cell.body.getRange().insertText(newValues[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Populate the \"Kit Components\" table with the correct reference/overview\n# values, replacing the placeholder content that had been written into the\n# wrong table, and blanking out the cells that no longer apply.\n$d = $word.ActiveDocument\n\n# Locate the \"Kit Components\" table (header row: Description / Quantity /\n# Volume / Storage of opened/reconstituted material) rather than assuming a\n# fixed index.\n$kitTable = $null\nforeach ($t in $d.Tables) {\n    if ($t.Rows.Item(1).Cells.Item(1).Range.Text.TrimEnd([char]7, [char]13) -eq \"Description\") {\n        $kitTable = $t\n        break\n    }\n}\n\n$rows = $kitTable.Rows\n\n$rows.Item(2).Cells.Item(1).Range.Text = \"Reactive Species\"\n$rows.Item(2).Cells.Item(2).Range.Text = \"Mouse\"\n$rows.Item(2).Cells.Item(3).Range.Text = \"\"\n$rows.Item(2).Cells.Item(4).Range.Text = \"\"\n\n$rows.Item(3).Cells.Item(1).Range.Text = \"Size\"\n$rows.Item(3).Cells.Item(2).Range.Text = \"96 wells/kit, with removable strips.\"\n$rows.Item(3).Cells.Item(3).Range.Text = \"\"\n$rows.Item(3).Cells.Item(4).Range.Text = \"\"\n\n$rows.Item(4).Cells.Item(1).Range.Text = \"Sensitivity*\"\n$rows.Item(4).Cells.Item(2).Range.Text = \"12 pg/ml\"\n$rows.Item(4).Cells.Item(3).Range.Text = \"\"\n$rows.Item(4).Cells.Item(4).Range.Text = \"\"\n\n$rows.Item(5).Cells.Item(1).Range.Text = \"Storage Instructions\"\n$rows.Item(5).Cells.Item(2).Range.Text = \"Store at 4\u00b0C for 6 months, at -20\u00b0C for 12 months. Avoid multiple freeze-thaw cycles (Ships with gel ice, can store for up to 3 days in room temperature. Freeze upon receiving.)\"\n$rows.Item(5).Cells.Item(3).Range.Text = \"\"\n$rows.Item(5).Cells.Item(4).Range.Text = \"\"\n\n$rows.Item(6).Cells.Item(1).Range.Text = \"Uniprot ID\"\n$rows.Item(6).Cells.Item(2).Range.Text = \"P15947\"\n$rows.Item(6).Cells.Item(3).Range.Text = \"\"\n$rows.Item(6).Cells.Item(4).Range.Text = \"\"\n\n$rows.Item(7).Cells.Item(1).Range.Text = \"\"\n$rows.Item(7).Cells.Item(2).Range.Text = \"\"\n$rows.Item(7).Cells.Item(3).Range.Text = \"\"\n$rows.Item(7).Cells.Item(4).Range.Text = \"\"\n\n$rows.Item(8).Cells.Item(1).Range.Text = \"\"\n$rows.Item(8).Cells.Item(2).Range.Text = \"\"\n$rows.Item(8).Cells.Item(3).Range.Text = \"\"\n$rows.Item(8).Cells.Item(4).Range.Text = \"\"\n"}
